# Bitstrip maker: redraw the "bitmap" grid (A1:I9), teach the T/E formulas
# to also pass through non-"X" marks (e.g. "P"), and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Redraw the bitmap grid (only the cells that actually changed) ---

# Row 3: C3 X, E3 cleared, G3 X
$ws.Range("C3").Value = "X"
$ws.Range("E3").Value = ""
$ws.Range("G3").Value = "X"

# Row 4: D4, E4, F4 cleared
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""

# Row 5: C5, D5 cleared; E5 becomes "P"; F5, G5 cleared
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "P"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""

# Row 6: D6, E6, F6 cleared
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

# Row 7: C7 X, E7 cleared, G7 X
$ws.Range("C7").Value = "X"
$ws.Range("E7").Value = ""
$ws.Range("G7").Value = "X"

# --- 2. Update the T/E formulas so a non-"X" mark passes through verbatim ---
# (was: IF(ISBLANK(ref), "E", "T")  ->  IF(ISBLANK(ref), "E", IF(ref="X", "T", ref)))

$ws.Range("A11").Formula = '=IF(ISBLANK(A1), "E", IF(A1="X", "T", A1))'
$ws.Range("B11:I11").Formula = '=IF(ISBLANK(B1), "E", IF(B1="X", "T", B1))'

$ws.Range("A12:I12").Formula = '=IF(ISBLANK(A2), "E", IF(A2="X", "T", A2))'
$ws.Range("A13:I13").Formula = '=IF(ISBLANK(A3), "E", IF(A3="X", "T", A3))'
$ws.Range("A14:I14").Formula = '=IF(ISBLANK(A4), "E", IF(A4="X", "T", A4))'
$ws.Range("A15:I15").Formula = '=IF(ISBLANK(A5), "E", IF(A5="X", "T", A5))'
$ws.Range("A16:I16").Formula = '=IF(ISBLANK(A6), "E", IF(A6="X", "T", A6))'
$ws.Range("A17:I17").Formula = '=IF(ISBLANK(A7), "E", IF(A7="X", "T", A7))'
$ws.Range("A18:I18").Formula = '=IF(ISBLANK(A8), "E", IF(A8="X", "T", A8))'
$ws.Range("A19:I19").Formula = '=IF(ISBLANK(A9), "E", IF(A9="X", "T", A9))'

# --- 3. Move the view: scroll down and park the selection on Q27 ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q27").Select()
